$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new parish row (row 7) with the "Pomy - Gressy - Suchy" data
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 7070000000
$ws.Range("C7").Value = 7070
$ws.Range("D7").Value = 7000
$ws.Range("E7").Value = "Pomy – Gressy – Suchy"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = "P"

# Move the active cell selection from G7 to E7
$ws.Range("E7").Select()
